# Generate Report for Handoff
#
# A new handoff run produced a fresh source-file GUID and a fresh content
# hash for the localization XLIFF packages, plus updated handoff
# timestamps for the zh-cn and de-de locales. Refresh every cell (and the
# matching hyperlink display text) that references the old identifiers.

$wb = $excel.ActiveWorkbook

$oldGuid = "71297eed-2fcd-4667-9fbd-e49f71006195"
$newGuid = "9c0ae9dd-fdc5-435a-9f7c-239370f426dd"

$oldHash = "c857f7687623ae3e7aed4dbdbc6c85559a3db1ee"
$newHash = "81fc26cf45b7e48f86bd683342bedf037d512b8f"

$newMd = "$newGuid.md"

$newFileZh = "$newGuid.$newHash.zh-cn.xlf"
$newFileDe = "$newGuid.$newHash.de-de.xlf"

$newTimeZh = "2016-03-09 10:53:52"
$newTimeDe = "2016-03-09 10:53:57"

# --- Overview sheet: A2 is the hyperlinked source .md filename ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newMd

# --- zh-cn sheet: A2 source filename, C2 handoff file, D2 handoff datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("C2").Value = $newFileZh
$wsZh.Range("D2").Value = $newTimeZh
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsZh.Hyperlinks.Item(2).TextToDisplay = $newFileZh

# --- de-de sheet: A2 source filename, C2 handoff file, D2 handoff datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("C2").Value = $newFileDe
$wsDe.Range("D2").Value = $newTimeDe
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsDe.Hyperlinks.Item(2).TextToDisplay = $newFileDe
